$d = $word.ActiveDocument

# Rewrite the first paragraph (currently a single run "Texto jajaj") as two
# runs -- "Texto " and "jajaj" -- with proofing-error bookmarks (spell-check
# squiggle markers) bracketing "jajaj", then append two blank paragraphs and
# a new paragraph with the added sentence. Using InsertXML lets us emit the
# exact WordprocessingML (including <w:proofErr/>) that Word's editor
# produces for this kind of edit.

$w_ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = '<w:p ' + $w_ns + '>' +
            '<w:r><w:t xml:space="preserve">Texto </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>jajaj</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
          '</w:p>' +
          '<w:p ' + $w_ns + '/>' +
          '<w:p ' + $w_ns + '/>' +
          '<w:p ' + $w_ns + '>' +
            '<w:r><w:t>Izan ha estado aquí</w:t></w:r>' +
          '</w:p>'

$target = $d.Paragraphs(1).Range
$target.InsertXML($newXml)
